$d = $word.ActiveDocument

# The resume is missing the contact-info line beneath the name header. Re-insert it as
# its own centered paragraph, directly after the "Dheeraj Chand" title paragraph (and
# before "PROFESSIONAL SUMMARY"), matching the long-resume layout.
$d.Content.Find.Execute(
    "Dheeraj Chand",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
